$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'66.313.97"
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.22%  '
$c = $ws.Range('D3')
$c.Value = "'3.205.22"
$c.Style = 'Normal'
$c = $ws.Range('D4')
$c.Value = "'1.00"
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$c = $ws.Range('D5')
$c.Value = "'608.38"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.82%  '
$c = $ws.Range('D6')
$c.Value = "'156.24"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.85%  '
$ws.Range('E7').Value = '  +0.04%  '
$c = $ws.Range('D8')
$c.Value = "'3.206.62"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.54%  '
$c = $ws.Range('D9')
$c.Value = "'0.550"
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -1.46%  '
$c = $ws.Range('D10')
$c.Value = "'0.161"
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +0.04%  '
$c = $ws.Range('D11')
$c.Value = "'5.65"
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -4.23%  '
$ws.Range('E12').Value = '  -3.08%  '
$ws.Range('E13').Value = '  -0.38%  '
$c = $ws.Range('D14')
$c.Value = "'38.45"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -2.21%  '
$c = $ws.Range('D15')
$c.Value = "'3.733.85"
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +0.54%  '
$c = $ws.Range('D16')
$c.Value = "'66.451.73"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('E17').Value = '  -2.81%  '
$c = $ws.Range('D18')
$c.Value = "'3.206.58"
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('E19').Value = '  +1.34%  '
$c = $ws.Range('D20')
$c.Value = "'506.43"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -2.63%  '
$c = $ws.Range('D21')
$c.Value = "'15.32"
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.89%  '
$c = $ws.Range('D22')
$c.Value = "'0.730"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.30%  '
$c = $ws.Range('D23')
$c.Value = "'7.99"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -1.69%  '
$ws.Range('E24').Value = '  -2.15%  '
$c = $ws.Range('D25')
$c.Value = "'85.18"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('E28').Value = '  -2.75%  '
$c = $ws.Range('D29')
$c.Value = "'2.35"
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -0.20%  '
$c = $ws.Range('D30')
$c.Value = "'0.128"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +41.44%  '
$c = $ws.Range('B31')
$c.Value = "'Stacks"
$c.Style = 'Normal'
$c = $ws.Range('C31')
$c.Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.Value = "'2.92"
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.81%  '
$c = $ws.Range('B32')
$c.Value = "'NEARProtocol"
$c.Style = 'Normal'
$c = $ws.Range('C32')
$c.Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.Value = "'6.95"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -1.50%  '
$c = $ws.Range('D33')
$c.Value = "'28.25"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  -5.26%  '
$c = $ws.Range('D36')
$c.Value = "'6.44"
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -1.67%  '
$c = $ws.Range('D37')
$c.Value = "'501.92"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -1.89%  '
$ws.Range('E38').Value = '  +0.98%  '
$c = $ws.Range('D39')
$c.Value = "'0.0₃0768"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +13.24%  '
$ws.Range('E40').Value = '  +1.97%  '
$c = $ws.Range('D41')
$c.Value = "'0.0419"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -1.71%  '
$ws.Range('E42').Value = '  +5.16%  '
$c = $ws.Range('D43')
$c.Value = "'8.71"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -2.45%  '
$ws.Range('E44').Value = '  -1.95%  '
$c = $ws.Range('D45')
$c.Value = "'2.913.25"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('E46').Value = '  -1.28%  '
$c = $ws.Range('D47')
$c.Value = "'28.17"
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.35%  '
$c = $ws.Range('D48')
$c.Value = "'2.40"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +2.16%  '
$ws.Range('E50').Value = '  -0.98%  '
$c = $ws.Range('D51')
$c.Value = "'122.15"
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.28%  '
